# Add lab test grades for students who were previously missing R (lab test)
# scores, and correct a couple of existing entries that now cross the
# pass/fail threshold. T/U columns recalc automatically for rows that carry
# the (R+S)/2+Q formula (rows 2-3 and 8-23); rows 4-7 never had that formula
# (it was only ever set up starting at row 8), so their Lab test mark (T)
# and Situation (U) are entered directly, matching the plain numbers/values
# that the workbook author typed in by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: existing grade correction - formula-driven recalculation
$ws.Range("R2").Value = 10

# Row 4: brand-new grade, no live formula on this row -> literal T/U
$ws.Range("R4").Value = 1
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = "fail"

# Row 5
$ws.Range("R5").Value = 8.8
$ws.Range("T5").Value = 8.8
$ws.Range("U5").Value = "pass"

# Row 6
$ws.Range("R6").Value = 6.5
$ws.Range("T6").Value = 6.5
$ws.Range("U6").Value = "pass"

# Row 7
$ws.Range("R7").Value = 6.4
$ws.Range("T7").Value = 6.4
$ws.Range("U7").Value = "pass"

# Row 9
$ws.Range("R9").Value = 5.6
$ws.Range("T9").Value = 5.6
$ws.Range("U9").Value = "pass"

# Row 11
$ws.Range("R11").Value = 6.3
$ws.Range("T11").Value = 6.3
$ws.Range("U11").Value = "pass"

# Row 13
$ws.Range("R13").Value = 7.1
$ws.Range("T13").Value = 7.1
$ws.Range("U13").Value = "pass"

# Row 14
$ws.Range("R14").Value = 7.1
$ws.Range("T14").Value = 7.1
$ws.Range("U14").Value = "pass"

# Row 16
$ws.Range("R16").Value = 5.9
$ws.Range("T16").Value = 5.9
$ws.Range("U16").Value = "pass"

# Row 18: existing grade correction - formula-driven recalculation (shared formula)
$ws.Range("S18").Value = 9.5

# Row 19: existing grade correction - formula-driven recalculation (shared formula)
$ws.Range("R19").Value = 8

# Row 20: brand-new grade, no live formula on this row -> literal T/U
$ws.Range("R20").Value = 1
$ws.Range("T20").Value = 1
$ws.Range("U20").Value = "fail"

# Row 21
$ws.Range("R21").Value = 6.5
$ws.Range("T21").Value = 7
$ws.Range("U21").Value = "pass"

# Row 22
$ws.Range("R22").Value = 7.1
$ws.Range("T22").Value = 7.1
$ws.Range("U22").Value = "pass"

# Match the author's final selection state (R2:R23, active cell R2)
$ws.Range("R2:R23").Select()
